$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume refresh (GitHub Actions data pull).
# D-column cells whose new value is a plain decimal number (e.g. "397.97") are written
# with a leading apostrophe so Excel stores them as literal TEXT (matching the sheet's
# existing inline-string convention) instead of auto-converting to a numeric value -
# then the cell Style is reset to "Normal" so no stray NumberFormat/quote-prefix style
# is left behind on the cell.

$ws.Range("D2").Value = '57.123.94'
$ws.Range("E2").Value = '  +1.75%  '
$ws.Range("D3").Value = '3.257.67'
$ws.Range("E3").Value = '  +1.24%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''397.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").Value = '''108.65'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.48%  '
$ws.Range("D7").Value = '''0.578'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.20%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -0.79%  '
$ws.Range("D10").Value = '''39.33'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.23%  '
$ws.Range("D11").Value = '''0.0954'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.01%  '
$ws.Range("E12").Value = '  +1.64%  '
$ws.Range("D13").Value = '3.773.16'
$ws.Range("E13").Value = '  +1.26%  '
$ws.Range("D14").Value = '''8.26'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.43%  '
$ws.Range("D15").Value = '''18.97'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.69%  '
$ws.Range("D16").Value = '3.257.52'
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("E17").Value = '  -2.30%  '
$ws.Range("D18").Value = '''11.04'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.01%  '
$ws.Range("D19").Value = '56.925.66'
$ws.Range("E19").Value = '  +1.78%  '
$ws.Range("D20").Value = '''3.31'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.56%  '
$ws.Range("E21").Value = '  +4.99%  '
$ws.Range("E22").Value = '  -1.22%  '
$ws.Range("D23").Value = '''293.43'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.24%  '
$ws.Range("D24").Value = '''74.06'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.80%  '
$ws.Range("D25").Value = '''3.18'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.72%  '
$ws.Range("D26").Value = '''7.94'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.38%  '
$ws.Range("D27").Value = '''28.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = '  +0.76%  '
$ws.Range("D29").Value = '''7.48'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.43%  '
$ws.Range("E30").Value = '  -2.78%  '
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("E32").Value = '  +1.02%  '
$ws.Range("E33").Value = '  -0.57%  '
$ws.Range("E34").Value = '  +10.63%  '
$ws.Range("E35").Value = '  -0.34%  '
$ws.Range("E36").Value = '  +0.83%  '
$ws.Range("D37").Value = '''51.28'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.54%  '
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("D39").Value = '''3.46'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.00%  '
$ws.Range("D40").Value = '''2.99'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.29%  '
$ws.Range("D41").Value = '''137.10'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.87%  '
$ws.Range("E42").Value = '  +1.55%  '
$ws.Range("E43").Value = '  -0.97%  '
$ws.Range("E44").Value = '  -2.92%  '
$ws.Range("D45").Value = '''3.92'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.54%  '
$ws.Range("D46").Value = '''16.78'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.44%  '
$ws.Range("D47").Value = '''22.38'
$ws.Range("D47").Style = "Normal"
$ws.Range("E48").Value = '  +4.90%  '
$ws.Range("D49").Value = '2.144.28'
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("D50").Value = '''2.44'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.85%  '
$ws.Range("D51").Value = '''1.98'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.00%  '
